$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowRange {
    param($ws, $row1, $row2, $colStart, $colEnd)
    $range1 = $ws.Range($colStart + $row1 + ":" + $colEnd + $row1)
    $range2 = $ws.Range($colStart + $row2 + ":" + $colEnd + $row2)
    $vals1 = $range1.Value2
    $vals2 = $range2.Value2
    $range1.Value2 = $vals2
    $range2.Value2 = $vals1
}

# The rows' match-details (home/away teams, scores, odds, timestamps, url) were
# re-sequenced: swap the F:V payload between the following row pairs while
# leaving the Indice/pais/torneio/temporada/data_partida (A:E) columns intact.
Swap-RowRange $ws 104 106 "F" "V"
Swap-RowRange $ws 113 114 "F" "V"
Swap-RowRange $ws 118 119 "F" "V"

# Append the new match row (index 120) after the former last row (120),
# copying row 120's formatting (bold/bordered index cell, date format) then
# overwriting the values for the new fixture.
$ws.Range("A120:V120").Copy($ws.Range("A121:V121"))

$ws.Range("A121").Value2 = 120
$ws.Range("B121").Value2 = "italy"
$ws.Range("C121").Value2 = "serie-c-group-c"
$ws.Range("D121").Value2 = "2023-2024"
$ws.Range("E121").Value2 = 45236.86458333334
$ws.Range("F121").Value2 = "Monopoli"
$ws.Range("G121").Value2 = 2
$ws.Range("H121").Value2 = "Audace Cerignola"
$ws.Range("I121").Value2 = 2
$ws.Range("J121").Value2 = 2.58
$ws.Range("K121").Value2 = "02/11/2023 08:12"
$ws.Range("L121").Value2 = 2.88
$ws.Range("M121").Value2 = "06/11/2023 20:41"
$ws.Range("N121").Value2 = 2.79
$ws.Range("O121").Value2 = "02/11/2023 08:12"
$ws.Range("P121").Value2 = 2.83
$ws.Range("Q121").Value2 = "06/11/2023 20:41"
$ws.Range("R121").Value2 = 2.82
$ws.Range("S121").Value2 = "02/11/2023 08:12"
$ws.Range("T121").Value2 = 2.8
$ws.Range("U121").Value2 = "06/11/2023 20:41"
$ws.Range("V121").Value2 = "https://www.betexplorer.com/football/italy/serie-c-group-c/monopoli-audace-cerignola/nN1VHtXj/"
